$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Source table")
$fm     = $wb.Worksheets.Item("FM mkdocs table")

# --- "Source table" sheet: the actual source-of-truth data changes ---

# Row 48: Observation station file (new) -> now fully read/write supported (X),
# added as of 0.3.0, backed by ObservationPointModel (obsFile).
$source.Range("B48").Value = "X"
$source.Range("C48").Value = "X"
$source.Range("D48").Value = "0.3.0"
$source.Range("E48").Value = "hydrolib.core.io.obs.models"
$source.Range("F48").Value = "ObservationPointModel"

# Row 59: RainfallRunoffModel moved from hydrolib.core.io.fnm.models to
# hydrolib.core.io.rr.models.
$source.Range("E59").Value = "hydrolib.core.io.rr.models"
$source.Range("G59").Value = "Used to be in hydrolib.core.io.fnm.models before 0.3.0"

# Row 60: BuiModel moved from hydrolib.core.io.bui.models to
# hydrolib.core.io.rr.meteo.models.
$source.Range("E60").Value = "hydrolib.core.io.rr.meteo.models"
$source.Range("G60").Value = "Used to be in hydrolib.core.io.bui.models before 0.3.0"

# --- Selection / active-cell bookkeeping to match the saved view state ---
$source.Range("B29").Select()
$fm.Range("A53").Select()

Write-Host "done"
